# Remove whitespace in the GapType column ("Arbeit, Privat" -> "Arbeit,Privat")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QAGlist_Teil1")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# GapType column is "M" (Gap1_type)
$col = "M"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("$col$r")
    if ($cell.Value() -eq "Arbeit, Privat") {
        $cell.Value = "Arbeit,Privat"
    }
}
